$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("paper")

$ws.Range("B2").Value = 3.6770716
$ws.Range("C2").Value = 10.4493092

$ws.Range("B3").Value = 3.8431654
$ws.Range("C3").Value = 21.0865803

$ws.Range("B4").Value = 7.8961658
$ws.Range("C4").Value = 24.9691487

$ws.Range("B5").Value = 2.7493126
$ws.Range("C5").Value = 3.9613209
$ws.Range("D5").ClearContents()

$ws.Range("B6").Value = 3.4504604
$ws.Range("C6").Value = 6.6732969

$ws.Range("C7").Value = 4.3620868

$ws.Range("C8").Value = 3.6701993

$ws.Range("C9").Value = 5.4564957

$ws.Range("C10").Value = 6.5406227

$ws.Range("B11").Value = 2.1554668
$ws.Range("C11").Value = 1.9343046

$ws.Range("B12").Value = 3.2669837
$ws.Range("C12").Value = 5.8615834

$ws.Range("B13").Value = -1.5730766
$ws.Range("C13").Value = -0.9853803

$ws.Range("B14").Value = 18.3928806
$ws.Range("C14").Value = 16.0873723

$ws.Range("B15").Value = 3.8713396
$ws.Range("C15").Value = 8.7933913

$ws.Range("B16").Value = 3.3037158
$ws.Range("C16").Value = 5.7348456

$ws.Range("B17").Value = 3.2516518
$ws.Range("C17").Value = 10.2305804
$ws.Range("D17").ClearContents()

$ws.Range("B18").Value = 2.8585616
$ws.Range("C18").Value = 6.4040938

$ws.Range("B19").Value = 5.1163425
$ws.Range("C19").Value = 17.7594201

$ws.Range("B20").Value = 5.5030878
$ws.Range("C20").Value = 8.7002999
$ws.Range("D20").ClearContents()

$ws.Range("B21").Value = 3.0226693
$ws.Range("C21").Value = 7.4471963

$ws.Range("B22").Value = 4.9562555
$ws.Range("C22").Value = 20.1476858

$ws.Range("B23").Value = 6.2428234
$ws.Range("C23").Value = 17.6932619

$ws.Range("C24").Value = 14.5233804

$ws.Range("B25").Value = 14.0802571
$ws.Range("C25").Value = 11.7756491
$ws.Range("D25").ClearContents()

$ws.Range("B26").Value = 11.4523993
$ws.Range("C26").Value = 6.4765879

$ws.Range("B29").Value = 0.6250278
$ws.Range("C29").Value = 2.4420417

$ws.Range("B31").Value = 3.8198693
$ws.Range("C31").Value = 6.9261068

$ws.Range("B34").Value = 2.6851441
$ws.Range("C34").Value = -6.8902

$ws.Range("B35").Value = 21.0227329
$ws.Range("C35").Value = 13.0235623
